$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing quiz labels with topic names (rows 2-5)
$ws.Range("A2").Value = "Quiz1 (C-basics)"
$ws.Range("A3").Value = "Quiz2 (Arrays & Strings)"
$ws.Range("A4").Value = "Quiz3 (Conditions & loops)"
$ws.Range("A5").Value = "Quiz4 (Functions)"

# Add new rows 7-10 for Quiz5 - Quiz8
$ws.Range("A7").Value = "Quiz5 (Struct & Union & Enum)"
$ws.Range("B7").Value = 26
$ws.Range("C7").Value = 30

$ws.Range("A8").Value = "Quiz6 (Preprocecssor directives)"
$ws.Range("B8").Value = 34
$ws.Range("C8").Value = 39

$ws.Range("A9").Value = "Quiz7 (Pointer)"
$ws.Range("B9").Value = 59
$ws.Range("C9").Value = 69

$ws.Range("A10").Value = "Quiz8 (Embedded-C)"
$ws.Range("B10").Value = 12
$ws.Range("C10").Value = 16

# Widen column A to fit the longer labels (target stored width ~33.140625;
# the COM property quantizes to 1/6-character steps on export, so 32.3 is
# the input that lands closest to the target after that re-quantization)
$ws.Columns.Item(1).ColumnWidth = 32.3

# Update the selected cell
$ws.Range("B12").Select()
